$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Drop the stale "_GoBack" bookmark left on the first paragraph.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Append a new paragraph: "Simply testing github features"
#    ("github" is wrapped in proofErr spell-check markers, as Word
#    does for a word it doesn't recognise), followed by a trailing
#    empty paragraph.
# ------------------------------------------------------------------
$d.Content.InsertParagraphAfter()
$newPara = $d.Paragraphs($d.Paragraphs.Count).Range

$fragment = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/fragment.xml" pkg:contentType="application/xml"><pkg:xmlData><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Simply testing </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>github</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> features</w:t></w:r></w:p></pkg:xmlData></pkg:part></pkg:package>
'@
$newPara.InsertXML($fragment)
